$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the 5 data rows (2-6) with the newly imported CSV / category data ---
$ws.Range("A2").Value = 1269
$ws.Range("B2").Value = "B0005197"
$ws.Range("C2").Value = "Маніпулятори"
$ws.Range("D2").Value = "920-002643"
$ws.Range("E2").Value = "Logitech"
$ws.Range("F2").Value = "K120 Ukr"
$ws.Range("G2").Value = "Клавіатура Logitech K120 Ukr (920-002643)"
$ws.Range("H2").Value = "конструкція - мембранна, USB, англійська, українська, повнорозмірна, Клавіш - 104, вологостійкість, безшумне введення, Колір - чорний"
$ws.Range("I2").Value = 8.6
$ws.Range("K2").Value = "Клавіатури"
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 369
$ws.Range("N2").Value = "1"
$ws.Range("O2").Value = 24
$ws.Range("P2").Value = "1"
$ws.Range("R2").Value = "0"
$ws.Range("S2").Value = 158
$ws.Range("T2").Value = "https://opt.brain.com.ua/Klaviatura_Logitech_K120_920-002643-p47223.html"
$ws.Range("U2").Value = "8471606000"
$ws.Range("V2").Value = 1011
$ws.Range("W2").Value = 67
$ws.Range("X2").Value = "Клавіатура"
$ws.Range("Y2").Value = "3"
$ws.Range("Z2").Value = "Китай"
$ws.Range("AA2").Value = 369
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = "0"
$ws.Range("AD2").Value = "0"
$ws.Range("A3").Value = 7925
$ws.Range("B3").Value = "418"
$ws.Range("C3").Value = "Витратні матеріали оригінальні"
$ws.Range("D3").Value = "51645AE"
$ws.Range("E3").Value = "HP"
$ws.Range("F3").Value = "DJ No. 45 Black"
$ws.Range("G3").Value = "Картридж HP DJ No. 45 Black (51645AE)"
$ws.Range("H3").Value = "струменевий, оригінальний, Black, Сумісність - Hewlett Packard"
$ws.Range("I3").Value = 66
$ws.Range("K3").Value = "Картриджі"
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "0"
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = "1"
$ws.Range("R3").Value = "0"
$ws.Range("S3").Value = 218
$ws.Range("T3").Value = "https://opt.brain.com.ua/Kartridjh_HP_DJ_No_45_Black_51645AE-p19822.html"
$ws.Range("U3").Value = "8443999090"
$ws.Range("V3").Value = 1034
$ws.Range("W3").Value = 102
$ws.Range("X3").Value = "Картридж"
$ws.Range("Y3").Value = "1"
$ws.Range("Z3").Value = "Китай"
$ws.Range("AA3").Value = 2844
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = "0"
$ws.Range("AD3").Value = "0"
$ws.Range("A4").Value = 7925
$ws.Range("B4").Value = "S0007198"
$ws.Range("C4").Value = "Витратні матеріали оригінальні"
$ws.Range("D4").Value = "CH561HE"
$ws.Range("E4").Value = "HP"
$ws.Range("F4").Value = "DJ No.122 Black, DJ 2050"
$ws.Range("G4").Value = "Картридж HP DJ No.122 Black, DJ 2050 (CH561HE)"
$ws.Range("H4").Value = "струменевий, оригінальний, Black, Сумісність - Hewlett Packard, 120 стр"
$ws.Range("I4").Value = 18.39
$ws.Range("K4").Value = "Картриджі"
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "0"
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = "1"
$ws.Range("R4").Value = "0"
$ws.Range("S4").Value = 239
$ws.Range("T4").Value = "https://opt.brain.com.ua/Kartridjh_HP_DJ_No122_Black_DJ_2050_CH561HE-p36184.html"
$ws.Range("U4").Value = "8443999090"
$ws.Range("V4").Value = 1034
$ws.Range("W4").Value = 102
$ws.Range("X4").Value = "Картридж"
$ws.Range("Y4").Value = "1"
$ws.Range("Z4").Value = "Китай"
$ws.Range("AA4").Value = 691
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = "0"
$ws.Range("AD4").Value = "0"
$ws.Range("A5").Value = 7925
$ws.Range("B5").Value = "S0007216"
$ws.Range("C5").Value = "Витратні матеріали оригінальні"
$ws.Range("D5").Value = "CH562HE"
$ws.Range("E5").Value = "HP"
$ws.Range("F5").Value = "DJ No.122 color, DJ 2050"
$ws.Range("G5").Value = "Картридж HP DJ No.122 color, DJ 2050 (CH562HE)"
$ws.Range("H5").Value = "струменевий, оригінальний, Color, Сумісність - Hewlett Packard, 100 стр"
$ws.Range("I5").Value = 20.2
$ws.Range("K5").Value = "Картриджі"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "0"
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = "1"
$ws.Range("R5").Value = "0"
$ws.Range("S5").Value = 240
$ws.Range("T5").Value = "https://opt.brain.com.ua/Kartridjh_HP_DJ_No122_color_DJ_2050_CH562HE-p36226.html"
$ws.Range("U5").Value = "3215902000"
$ws.Range("V5").Value = 1034
$ws.Range("W5").Value = 102
$ws.Range("X5").Value = "Картридж"
$ws.Range("Y5").Value = "1"
$ws.Range("Z5").Value = "Китай"
$ws.Range("AA5").Value = 798
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = "0"
$ws.Range("AD5").Value = "0"
$ws.Range("A6").Value = 7925
$ws.Range("B6").Value = "KM09159"
$ws.Range("C6").Value = "Витратні матеріали оригінальні"
$ws.Range("D6").Value = "2146B001/2146B005/21460001"
$ws.Range("E6").Value = "Canon"
$ws.Range("F6").Value = "CL-38 Color"
$ws.Range("G6").Value = "Картридж CL-38 Color Canon (2146B001/2146B005/21460001)"
$ws.Range("H6").Value = "струменевий, оригінальний, Magenta, Yellow, Cyan, Сумісність - Canon, 205 стр"
$ws.Range("I6").Value = 20.75
$ws.Range("K6").Value = "Картриджі"
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = "0"
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = "1"
$ws.Range("R6").Value = "0"
$ws.Range("S6").Value = 320
$ws.Range("T6").Value = "https://opt.brain.com.ua/Kartridjh_CANON_CL-38_Color_2146B001_2146B005_21460001-p19728.html"
$ws.Range("U6").Value = "8443999090"
$ws.Range("V6").Value = 1034
$ws.Range("W6").Value = 102
$ws.Range("X6").Value = "Картридж"
$ws.Range("Y6").Value = "2"
$ws.Range("Z6").Value = "Китай"
$ws.Range("AA6").Value = 795
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = "0"
$ws.Range("AD6").Value = "0"

# --- Selection moves to C12 (matches the author's last click before saving) ---
$ws.Range("C12").Select()

# --- Shrink the AutoFilter range down to the header row only (A1:AD1) ---
$ws.AutoFilterMode = $false
$ws.Range("A1:AD1").AutoFilter()

# --- The _xlnm._FilterDatabase defined name must track the new AutoFilter range ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=" + $ws.Name + "!`$A`$1:`$AD`$1"
    }
}
